$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 6: wake-up time (E6) and derived formulas (C6, F6)
$ws.Range("C6").Formula = "=(E6-D6)*1440"
$ws.Range("E6").Value = 0.3625
$ws.Range("F6").Formula = "=C6/B6"

# Update the selected cell in the sheet view to B19
$ws.Range("B19").Select()
